# Insert a new weekly record at row 32 ("Hortaliza, Terminal La Palmera de
# La Serena - Berenjena"), shifting the existing rows 32-102 down to 33-103.
# The new row duplicates the market/category metadata of the (old) row 37
# record but carries a new, later reporting date (serial 44526 = 2021-11-26).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 32..102 down to 33..103, leaving a blank row 32 to fill in.
$ws.Rows.Item(32).Insert()

$epoch = Get-Date -Year 1899 -Month 12 -Day 30 -Hour 0 -Minute 0 -Second 0
$fecha = $epoch.AddDays(44526)

$ws.Range("A32").Value = 8
$ws.Range("B32").Value = "Terminal La Palmera de La Serena"
$ws.Range("C32").Value = "Coquimbo"
$ws.Range("D32").Value = $fecha
$ws.Range("E32").Value = 4
$ws.Range("F32").Value = 100112001
$ws.Range("G32").Value = "Berenjena"
$ws.Range("H32").Value = "Sin especificar"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 600
$ws.Range("K32").Value = 9000
$ws.Range("L32").Value = 10000
$ws.Range("M32").Value = 9500
$ws.Range("N32").Value = "`$/caja 60 unidades"
$ws.Range("O32").Value = "Región de Arica y Parinacota"
$ws.Range("P32").Value = 158
$ws.Range("Q32").Value = 60
$ws.Range("R32").Value = "Hortaliza"
